$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 106
$ws_ALC.Range("H106").Value = 35295860
$ws_ALC.Range("I106").Value = 40001604
$ws_ALC.Range("J106").Value = 2800
$ws_ALC.Range("K106").Value = 40001604
$ws_ALC.Range("L106").Value = 2800
$ws_ALC.Range("M106").Value = -40000973
$ws_ALC.Range("N106").Value = -4062

# ALC row 113
$ws_ALC.Range("H113").Value = 2116.6667
$ws_ALC.Range("I113").Value = 2116.6667
$ws_ALC.Range("J113").Value = 0
$ws_ALC.Range("K113").Value = 2116.6667
$ws_ALC.Range("L113").Value = 0
$ws_ALC.Range("M113").Value = 1137.3333
$ws_ALC.Range("N113").ClearContents()

# ALC row 138
$ws_ALC.Range("H138").Value = 2552.8152
$ws_ALC.Range("I138").Value = 3115.077
$ws_ALC.Range("J138").Value = 2460.2913
$ws_ALC.Range("K138").Value = 9345.231
$ws_ALC.Range("L138").Value = 7380.8739
$ws_ALC.Range("M138").Value = -4205.231
$ws_ALC.Range("N138").Value = -17660.8739

# ALC row 141
$ws_ALC.Range("H141").Value = 6840.4736
$ws_ALC.Range("I141").Value = 3031.5833
$ws_ALC.Range("K141").Value = 9094.749899999999
$ws_ALC.Range("M141").Value = -3914.749899999999

# ARM row 2
$ws_ARM.Range("H2").Value = 2386.3845
$ws_ARM.Range("I2").Value = 2502.3
$ws_ARM.Range("J2").Value = 2000
$ws_ARM.Range("K2").Value = 2502.3
$ws_ARM.Range("L2").Value = 2000
$ws_ARM.Range("M2").Value = -2389.3
$ws_ARM.Range("N2").Value = -2226

# ARM row 61
$ws_ARM.Range("H61").Value = 3200.7036
$ws_ARM.Range("I61").Value = 2772.7334
$ws_ARM.Range("J61").Value = 3735.6667
$ws_ARM.Range("K61").Value = 2772.7334
$ws_ARM.Range("L61").Value = 3735.6667
$ws_ARM.Range("M61").Value = -2560.7334
$ws_ARM.Range("N61").Value = -4159.6667

# ARM row 74
$ws_ARM.Range("H74").Value = 2627.4285
$ws_ARM.Range("I74").Value = 2482.9473
$ws_ARM.Range("J74").Value = 4000
$ws_ARM.Range("K74").Value = 2482.9473
$ws_ARM.Range("L74").Value = 4000
$ws_ARM.Range("M74").Value = -1608.9473
$ws_ARM.Range("N74").Value = -5748

# ARM row 77
$ws_ARM.Range("H77").Value = 2627.4285
$ws_ARM.Range("I77").Value = 2482.9473
$ws_ARM.Range("J77").Value = 4000
$ws_ARM.Range("K77").Value = 12414.7365
$ws_ARM.Range("L77").Value = 20000
$ws_ARM.Range("M77").Value = -8046.736499999999
$ws_ARM.Range("N77").Value = -28736

# ARM row 97
$ws_ARM.Range("H97").Value = 789.7368
$ws_ARM.Range("I97").Value = 694.41174
$ws_ARM.Range("J97").Value = 1600
$ws_ARM.Range("K97").Value = 694.41174
$ws_ARM.Range("L97").Value = 1600
$ws_ARM.Range("M97").Value = -198.41174
$ws_ARM.Range("N97").Value = -2592

# ARM row 102
$ws_ARM.Range("H102").Value = 2056.6316
$ws_ARM.Range("I102").Value = 1980.375
$ws_ARM.Range("J102").Value = 2463.3333
$ws_ARM.Range("K102").Value = 1980.375
$ws_ARM.Range("L102").Value = 2463.3333
$ws_ARM.Range("M102").Value = -358.375
$ws_ARM.Range("N102").Value = -5707.3333

# ARM row 110
$ws_ARM.Range("H110").Value = 1367.5
$ws_ARM.Range("I110").Value = 1367.5
$ws_ARM.Range("K110").Value = 1367.5
$ws_ARM.Range("M110").Value = 677.5

# ARM row 116
$ws_ARM.Range("H116").Value = 2386.3845
$ws_ARM.Range("I116").Value = 2502.3
$ws_ARM.Range("J116").Value = 2000
$ws_ARM.Range("K116").Value = 2502.3
$ws_ARM.Range("L116").Value = 2000
$ws_ARM.Range("M116").Value = -208.3000000000002
$ws_ARM.Range("N116").Value = -6588

# ARM row 136
$ws_ARM.Range("H136").Value = 3200.7036
$ws_ARM.Range("I136").Value = 2772.7334
$ws_ARM.Range("J136").Value = 3735.6667
$ws_ARM.Range("K136").Value = 8318.200199999999
$ws_ARM.Range("L136").Value = 11207.0001
$ws_ARM.Range("M136").Value = -5768.200199999999
$ws_ARM.Range("N136").Value = -16307.0001

# BSM row 3
$ws_BSM.Range("H3").Value = 2386.3845
$ws_BSM.Range("I3").Value = 2502.3
$ws_BSM.Range("J3").Value = 2000
$ws_BSM.Range("K3").Value = 2502.3
$ws_BSM.Range("L3").Value = 2000
$ws_BSM.Range("M3").Value = -2388.3
$ws_BSM.Range("N3").Value = -2228

# BSM row 99
$ws_BSM.Range("H99").Value = 948.5909
$ws_BSM.Range("I99").Value = 842.36365
$ws_BSM.Range("J99").Value = 1054.8182
$ws_BSM.Range("K99").Value = 842.36365
$ws_BSM.Range("L99").Value = 1054.8182
$ws_BSM.Range("M99").Value = 655.63635
$ws_BSM.Range("N99").Value = -4050.8182

# BSM row 105
$ws_BSM.Range("H105").Value = 12503276
$ws_BSM.Range("I105").Value = 17860310
$ws_BSM.Range("J105").Value = 3530
$ws_BSM.Range("K105").Value = 17860310
$ws_BSM.Range("L105").Value = 3530
$ws_BSM.Range("M105").Value = -17858563
$ws_BSM.Range("N105").Value = -7024

# BSM row 107
$ws_BSM.Range("H107").Value = 46579.09
$ws_BSM.Range("I107").Value = 56629.332
$ws_BSM.Range("J107").Value = 1353
$ws_BSM.Range("K107").Value = 56629.332
$ws_BSM.Range("L107").Value = 1353
$ws_BSM.Range("M107").Value = -54709.332
$ws_BSM.Range("N107").Value = -5193

# CRP row 7
$ws_CRP.Range("H7").Value = 75.166664
$ws_CRP.Range("I7").Value = 55.625
$ws_CRP.Range("J7").Value = 114.25
$ws_CRP.Range("K7").Value = 55.625
$ws_CRP.Range("L7").Value = 114.25
$ws_CRP.Range("M7").Value = 57.375
$ws_CRP.Range("N7").Value = -340.25

# CRP row 22
$ws_CRP.Range("H22").Value = 370.5
$ws_CRP.Range("I22").Value = 287.85715
$ws_CRP.Range("J22").Value = 563.3333
$ws_CRP.Range("K22").Value = 287.85715
$ws_CRP.Range("L22").Value = 563.3333
$ws_CRP.Range("M22").Value = 62.14285000000001
$ws_CRP.Range("N22").Value = -1263.3333

# CRP row 31
$ws_CRP.Range("H31").Value = 6132.8477
$ws_CRP.Range("I31").Value = 1272.7916
$ws_CRP.Range("J31").Value = 11434.728
$ws_CRP.Range("K31").Value = 1272.7916
$ws_CRP.Range("L31").Value = 11434.728
$ws_CRP.Range("M31").Value = -977.7916
$ws_CRP.Range("N31").Value = -12024.728

# CRP row 34
$ws_CRP.Range("H34").Value = 6132.8477
$ws_CRP.Range("I34").Value = 1272.7916
$ws_CRP.Range("J34").Value = 11434.728
$ws_CRP.Range("K34").Value = 1272.7916
$ws_CRP.Range("L34").Value = 11434.728
$ws_CRP.Range("M34").Value = -1070.7916
$ws_CRP.Range("N34").Value = -11838.728

# CRP row 36
$ws_CRP.Range("H36").Value = 4448
$ws_CRP.Range("I36").Value = 4448
$ws_CRP.Range("K36").Value = 4448
$ws_CRP.Range("M36").Value = -4060

# CRP row 40
$ws_CRP.Range("H40").Value = 4448
$ws_CRP.Range("I40").Value = 4448
$ws_CRP.Range("K40").Value = 4448
$ws_CRP.Range("M40").Value = -4288

# CRP row 62
$ws_CRP.Range("H62").Value = 4600
$ws_CRP.Range("I62").Value = 4542.857
$ws_CRP.Range("J62").Value = 5000
$ws_CRP.Range("K62").Value = 4542.857
$ws_CRP.Range("L62").Value = 5000
$ws_CRP.Range("M62").Value = -3918.857
$ws_CRP.Range("N62").Value = -6248

# CRP row 65
$ws_CRP.Range("H65").Value = 4600
$ws_CRP.Range("I65").Value = 4542.857
$ws_CRP.Range("J65").Value = 5000
$ws_CRP.Range("K65").Value = 22714.285
$ws_CRP.Range("L65").Value = 25000
$ws_CRP.Range("M65").Value = -19594.285
$ws_CRP.Range("N65").Value = -31240

# CRP row 105
$ws_CRP.Range("H105").Value = 1999.8
$ws_CRP.Range("I105").Value = 1999.6666
$ws_CRP.Range("J105").Value = 2000
$ws_CRP.Range("K105").Value = 1999.6666
$ws_CRP.Range("L105").Value = 2000
$ws_CRP.Range("M105").Value = -252.6666
$ws_CRP.Range("N105").Value = -5494

# CRP row 107
$ws_CRP.Range("H107").Value = 3125798.8
$ws_CRP.Range("I107").Value = 4464959
$ws_CRP.Range("J107").Value = 1091.6666
$ws_CRP.Range("K107").Value = 4464959
$ws_CRP.Range("L107").Value = 1091.6666
$ws_CRP.Range("M107").Value = -4463039
$ws_CRP.Range("N107").Value = -4931.6666

# CUL row 8
$ws_CUL.Range("H8").Value = 222.625
$ws_CUL.Range("I8").Value = 222.625
$ws_CUL.Range("K8").Value = 667.875
$ws_CUL.Range("M8").Value = -528.875

# CUL row 114
$ws_CUL.Range("H114").Value = 1188.2941
$ws_CUL.Range("J114").Value = 2642.5715
$ws_CUL.Range("L114").Value = 7927.7145
$ws_CUL.Range("N114").Value = -14435.7145

# CUL row 117
$ws_CUL.Range("H117").Value = 2561.5557
$ws_CUL.Range("I117").Value = 280
$ws_CUL.Range("J117").Value = 2846.75
$ws_CUL.Range("K117").Value = 840
$ws_CUL.Range("L117").Value = 8540.25
$ws_CUL.Range("M117").Value = 2602
$ws_CUL.Range("N117").Value = -15424.25

# CUL row 129
$ws_CUL.Range("H129").Value = 1743.0454
$ws_CUL.Range("J129").Value = 1867.85
$ws_CUL.Range("L129").Value = 5603.549999999999
$ws_CUL.Range("N129").Value = -15603.55

# CUL row 131
$ws_CUL.Range("H131").Value = 1076.4412
$ws_CUL.Range("I131").Value = 500
$ws_CUL.Range("J131").Value = 1093.909
$ws_CUL.Range("K131").Value = 1500
$ws_CUL.Range("L131").Value = 3281.727
$ws_CUL.Range("M131").Value = 3540
$ws_CUL.Range("N131").Value = -13361.727

# GSM row 2
$ws_GSM.Range("H2").Value = 81.181816
$ws_GSM.Range("I2").Value = 49.25
$ws_GSM.Range("J2").Value = 99.42856999999999
$ws_GSM.Range("K2").Value = 49.25
$ws_GSM.Range("L2").Value = 99.42856999999999
$ws_GSM.Range("M2").Value = 63.75
$ws_GSM.Range("N2").Value = -325.42857

# GSM row 70
$ws_GSM.Range("H70").Value = 5664.3784
$ws_GSM.Range("I70").Value = 5910.7
$ws_GSM.Range("J70").Value = 5374.5884
$ws_GSM.Range("K70").Value = 5910.7
$ws_GSM.Range("L70").Value = 5374.5884
$ws_GSM.Range("M70").Value = -5640.7
$ws_GSM.Range("N70").Value = -5914.5884

# GSM row 73
$ws_GSM.Range("H73").Value = 5664.3784
$ws_GSM.Range("I73").Value = 5910.7
$ws_GSM.Range("J73").Value = 5374.5884
$ws_GSM.Range("K73").Value = 5910.7
$ws_GSM.Range("L73").Value = 5374.5884
$ws_GSM.Range("M73").Value = -4974.7
$ws_GSM.Range("N73").Value = -7246.5884

# GSM row 102
$ws_GSM.Range("H102").Value = 1170.6666
$ws_GSM.Range("I102").Value = 1170.6666
$ws_GSM.Range("J102").Value = 0
$ws_GSM.Range("K102").Value = 1170.6666
$ws_GSM.Range("L102").Value = 0
$ws_GSM.Range("M102").Value = 451.3334
$ws_GSM.Range("N102").ClearContents()

# WVR row 62
$ws_WVR.Range("H62").Value = 152000
$ws_WVR.Range("J62").Value = 152000
$ws_WVR.Range("L62").Value = 152000
$ws_WVR.Range("N62").Value = -153248

# WVR row 65
$ws_WVR.Range("H65").Value = 152000
$ws_WVR.Range("J65").Value = 152000
$ws_WVR.Range("L65").Value = 760000
$ws_WVR.Range("N65").Value = -766240

# WVR row 107
$ws_WVR.Range("H107").Value = 495.33334
$ws_WVR.Range("J107").Value = 509.66666
$ws_WVR.Range("L107").Value = 1528.99998
$ws_WVR.Range("N107").Value = -5368.999980000001

# WVR row 122
$ws_WVR.Range("H122").Value = 3095.5557
$ws_WVR.Range("I122").Value = 2500
$ws_WVR.Range("J122").Value = 3265.7144
$ws_WVR.Range("K122").Value = 7500
$ws_WVR.Range("L122").Value = 9797.143199999999
$ws_WVR.Range("M122").Value = -5050
$ws_WVR.Range("N122").Value = -14697.1432

# WVR row 126
$ws_WVR.Range("H126").Value = 1898.4706
$ws_WVR.Range("I126").Value = 1487.6364
$ws_WVR.Range("J126").Value = 2651.6667
$ws_WVR.Range("K126").Value = 4462.9092
$ws_WVR.Range("L126").Value = 7955.000100000001
$ws_WVR.Range("M126").Value = -1992.9092
$ws_WVR.Range("N126").Value = -12895.0001

# WVR row 132
$ws_WVR.Range("H132").Value = 5379074.5
$ws_WVR.Range("I132").Value = 3272.6428
$ws_WVR.Range("J132").Value = 9806206
$ws_WVR.Range("K132").Value = 9817.928400000001
$ws_WVR.Range("L132").Value = 29418618
$ws_WVR.Range("M132").Value = -7287.928400000001
$ws_WVR.Range("N132").Value = -29423678
